$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove yellow "new data" highlight formatting from cells whose date value
# did not change this update (they revert to the plain "not-updated" date style).
# C3 already carries that plain date style, so copy its format onto the targets.
$ws.Range("C3").Copy()
$ws.Range("C17").PasteSpecial(-4122)
$ws.Range("C18").PasteSpecial(-4122)
$ws.Range("N22").PasteSpecial(-4122)
$ws.Range("N23").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Row 29: T5YIFR (5yr, 5yr Forward) ---
$ws.Range("N29").Value = 46042
$ws.Range("Q29").Value = 2.26
$ws.Range("R29").Value = 2.27
$ws.Range("S29").Value = 2.22
$ws.Range("T29").Value = 2.22
$ws.Range("U29").Value = 2.23

# --- Row 30: T10YIE (10yr TIPS) ---
$ws.Range("N30").Value = 46042
$ws.Range("Q30").Value = 2.33
$ws.Range("R30").Value = 2.33
$ws.Range("S30").Value = 2.29
$ws.Range("T30").Value = 2.29
$ws.Range("U30").Value = 2.3

# --- Row 47: DFF (FFR) ---
$ws.Range("N47").Value = 46041

# --- Row 48: DGS2 (2y UST) ---
$ws.Range("N48").Value = 46038
$ws.Range("Q48").Value = 3.59
$ws.Range("R48").Value = 3.56
$ws.Range("S48").Value = 3.51
$ws.Range("T48").Value = 3.53

# --- Row 49: DGS5 (5y UST) ---
$ws.Range("N49").Value = 46038
$ws.Range("Q49").Value = 3.82
$ws.Range("R49").Value = 3.77
$ws.Range("S49").Value = 3.72
$ws.Range("T49").Value = 3.75
$ws.Range("U49").Value = 3.77

# --- Row 50: DGS10 (10y UST) ---
$ws.Range("N50").Value = 46038
$ws.Range("Q50").Value = 4.24
$ws.Range("R50").Value = 4.17
$ws.Range("S50").Value = 4.15
$ws.Range("T50").Value = 4.18
$ws.Range("U50").Value = 4.19

# --- Row 52: DBAA (BAA) ---
$ws.Range("N52").Value = 46038
$ws.Range("Q52").Value = 5.87
$ws.Range("R52").Value = 5.82
$ws.Range("S52").Value = 5.83
$ws.Range("T52").Value = 5.87
$ws.Range("U52").Value = 5.89
